{"js": "// Clear all existing content from the document body, leaving a single\n// empty paragraph (mirrors the target OOXML: <w:body><w:p/><w:sectPr/></w:body>).\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// Delete every paragraph from last to first so index shifting doesn't\n// skip any paragraph. Word always keeps a final paragraph mark, so the\n// body ends up with a single empty paragraph, matching the target.\nfor (let i = paragraphs.items.length - 1; i >= 0; i--) {\n  paragraphs.items[i].delete();\n}\nawait context.sync();\n", "ps1": "# Clear all existing content from the document body, leaving a single\n# empty paragraph (mirrors the target OOXML: <w:body><w:p/><w:sectPr/></w:body>).\n$d = $word.ActiveDocument\n\n$count = $d.Paragraphs.Count\nfor ($i = $count; $i -ge 1; $i--) {\n    $d.Paragraphs.Item($i).Range.Delete()\n}\n"}
